$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 336151.7
$ws.Range("J17").Value = 336151.7
$ws.Range("L17").Value = 1008455.1
$ws.Range("N17").Value = -1008791.1
# Row 55
$ws.Range("H55").Value = 844.6
$ws.Range("I55").Value = 61.333332
$ws.Range("J55").Value = 1180.2858
$ws.Range("K55").Value = 61.333332
$ws.Range("L55").Value = 1180.2858
$ws.Range("M55").Value = 152.666668
$ws.Range("N55").Value = -1608.2858
# Row 93
$ws.Range("H93").Value = 75000
$ws.Range("J93").Value = 75000
$ws.Range("L93").Value = 75000
$ws.Range("N93").Value = -79992
# Row 95
$ws.Range("H95").Value = 50624
$ws.Range("J95").Value = 50624
$ws.Range("L95").Value = 50624
$ws.Range("N95").Value = -56116
# Row 103
$ws.Range("H103").Value = 668.5
$ws.Range("J103").Value = 933
$ws.Range("L103").Value = 2799
$ws.Range("N103").Value = -3971
# Row 106
$ws.Range("H106").Value = 2445
$ws.Range("I106").Value = 2333.9285
$ws.Range("K106").Value = 2333.9285
$ws.Range("M106").Value = -1702.9285
# Row 132
$ws.Range("H132").Value = 14930350
$ws.Range("I132").Value = 16951696
$ws.Range("J132").Value = 22921.75
$ws.Range("K132").Value = 50855088
$ws.Range("L132").Value = 68765.25
$ws.Range("M132").Value = -50852558
$ws.Range("N132").Value = -73825.25
# Row 138
$ws.Range("H138").Value = 1377625
$ws.Range("J138").Value = 1669833.4
$ws.Range("L138").Value = 5009500.199999999
$ws.Range("N138").Value = -5019780.199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17087.643
$ws.Range("I32").Value = 19021
$ws.Range("K32").Value = 19021
$ws.Range("M32").Value = -18734
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
# Row 132
$ws.Range("H132").Value = 7163.933
$ws.Range("I132").Value = 2121.75
$ws.Range("K132").Value = 6365.25
$ws.Range("M132").Value = -3835.25
# Row 135
$ws.Range("H135").Value = 66615
$ws.Range("J135").Value = 66615
$ws.Range("L135").Value = 66615
$ws.Range("N135").Value = -76755

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 95
$ws.Range("H95").Value = 62500
$ws.Range("J95").Value = 62500
$ws.Range("L95").Value = 62500
$ws.Range("N95").Value = -67992
# Row 96
$ws.Range("H96").Value = 39286.75
$ws.Range("I96").Value = 3574.5
$ws.Range("K96").Value = 3574.5
$ws.Range("M96").Value = -828.5
# Row 97
$ws.Range("H97").Value = 7663
$ws.Range("I97").Value = 1532.8334
$ws.Range("J97").Value = 44444
$ws.Range("K97").Value = 1532.8334
$ws.Range("L97").Value = 44444
$ws.Range("M97").Value = -541.8334
$ws.Range("N97").Value = -46426
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
# Row 105
$ws.Range("H105").Value = 2803.1875
$ws.Range("I105").Value = 2083.6667
$ws.Range("J105").Value = 6688.6
$ws.Range("K105").Value = 2083.6667
$ws.Range("L105").Value = 6688.6
$ws.Range("M105").Value = -336.6667000000002
$ws.Range("N105").Value = -10182.6
# Row 134
$ws.Range("H134").Value = 4671.027
$ws.Range("I134").Value = 1344.6279
$ws.Range("J134").Value = 9285.064
$ws.Range("K134").Value = 4033.8837
$ws.Range("L134").Value = 27855.192
$ws.Range("M134").Value = -1498.8837
$ws.Range("N134").Value = -32925.192
# Row 140
$ws.Range("H140").Value = 176319.75
$ws.Range("I140").Value = 74499
$ws.Range("J140").Value = 210260
$ws.Range("K140").Value = 74499
$ws.Range("L140").Value = 210260
$ws.Range("M140").Value = -69319
$ws.Range("N140").Value = -220620

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2574.75
$ws.Range("I31").Value = 1962.878
$ws.Range("J31").Value = 6158.5713
$ws.Range("K31").Value = 1962.878
$ws.Range("L31").Value = 6158.5713
$ws.Range("M31").Value = -1667.878
$ws.Range("N31").Value = -6748.5713
# Row 34
$ws.Range("H34").Value = 2574.75
$ws.Range("I34").Value = 1962.878
$ws.Range("J34").Value = 6158.5713
$ws.Range("K34").Value = 1962.878
$ws.Range("L34").Value = 6158.5713
$ws.Range("M34").Value = -1760.878
$ws.Range("N34").Value = -6562.5713
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = 0
# Row 51
$ws.Range("H51").Value = 48675
$ws.Range("J51").Value = 87350
$ws.Range("L51").Value = 87350
$ws.Range("N51").Value = -88822
# Row 58
$ws.Range("H58").Value = 2181.913
$ws.Range("I58").Value = 2214.3845
$ws.Range("J58").Value = 2139.7
$ws.Range("K58").Value = 2214.3845
$ws.Range("L58").Value = 2139.7
$ws.Range("M58").Value = -2011.3845
$ws.Range("N58").Value = -2545.7
# Row 61
$ws.Range("H61").Value = 48675
$ws.Range("J61").Value = 87350
$ws.Range("L61").Value = 87350
$ws.Range("N61").Value = -88046
# Row 132
$ws.Range("H132").Value = 1253171.5
$ws.Range("I132").Value = 1541799.8
$ws.Range("J132").Value = 2448.8333
$ws.Range("K132").Value = 4625399.4
$ws.Range("L132").Value = 7346.499899999999
$ws.Range("M132").Value = -4622869.4
$ws.Range("N132").Value = -12406.4999
# Row 134
$ws.Range("H134").Value = 3145.8167
$ws.Range("I134").Value = 1778.8667
$ws.Range("J134").Value = 7246.6665
$ws.Range("K134").Value = 5336.6001
$ws.Range("L134").Value = 21739.9995
$ws.Range("M134").Value = -2801.6001
$ws.Range("N134").Value = -26809.9995
# Row 136
$ws.Range("H136").Value = 2181.913
$ws.Range("I136").Value = 2214.3845
$ws.Range("J136").Value = 2139.7
$ws.Range("K136").Value = 6643.1535
$ws.Range("L136").Value = 6419.099999999999
$ws.Range("M136").Value = -4093.1535
$ws.Range("N136").Value = -11519.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 71678.07000000001
$ws.Range("I11").Value = 83545.664
$ws.Range("J11").Value = 472.5
$ws.Range("K11").Value = 250636.992
$ws.Range("L11").Value = 1417.5
$ws.Range("M11").Value = -250496.992
$ws.Range("N11").Value = -1697.5
# Row 55
$ws.Range("H55").Value = 149641.42
$ws.Range("J55").Value = 149641.42
$ws.Range("L55").Value = 448924.26
$ws.Range("N55").Value = -449278.26
# Row 86
$ws.Range("H86").Value = 621.0833
$ws.Range("J86").Value = 533
$ws.Range("L86").Value = 1599
$ws.Range("N86").Value = -3971
# Row 89
$ws.Range("H89").Value = 621.0833
$ws.Range("J89").Value = 533
$ws.Range("L89").Value = 4797
$ws.Range("N89").Value = -16653

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 94
$ws.Range("H94").Value = 56798.332
$ws.Range("J94").Value = 56798.332
$ws.Range("L94").Value = 56798.332
$ws.Range("N94").Value = -58150.332
# Row 132
$ws.Range("H132").Value = 3987.32
$ws.Range("I132").Value = 3838.5217
$ws.Range("J132").Value = 5698.5
$ws.Range("K132").Value = 11515.5651
$ws.Range("L132").Value = 17095.5
$ws.Range("M132").Value = -8985.5651
$ws.Range("N132").Value = -22155.5
# Row 133
$ws.Range("H133").Value = 76926.664
$ws.Range("J133").Value = 76926.664
$ws.Range("L133").Value = 76926.664
$ws.Range("N133").Value = -87046.664
# Row 140
$ws.Range("H140").Value = 86250
$ws.Range("J140").Value = 86250
$ws.Range("L140").Value = 86250
$ws.Range("N140").Value = -96610

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 131.11539
$ws.Range("I55").Value = 96.69231000000001
$ws.Range("K55").Value = 96.69231000000001
$ws.Range("M55").Value = 76.30768999999999
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
# Row 132
$ws.Range("H132").Value = 2471.3281
$ws.Range("I132").Value = 2439.3276
$ws.Range("J132").Value = 2780.6667
$ws.Range("K132").Value = 7317.9828
$ws.Range("L132").Value = 8342.000100000001
$ws.Range("M132").Value = -4787.9828
$ws.Range("N132").Value = -13402.0001
# Row 136
$ws.Range("H136").Value = 3664.2432
$ws.Range("I136").Value = 3197.724
$ws.Range("K136").Value = 9593.172
$ws.Range("M136").Value = -7043.172

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
# Row 97
$ws.Range("H97").Value = 69999
$ws.Range("J97").Value = 69999
$ws.Range("L97").Value = 69999
$ws.Range("N97").Value = -71981
# Row 100
$ws.Range("H100").Value = 1112.2439
$ws.Range("I100").Value = 1096.9697
$ws.Range("J100").Value = 1175.25
$ws.Range("K100").Value = 2193.9394
$ws.Range("L100").Value = 2350.5
$ws.Range("M100").Value = -1652.9394
$ws.Range("N100").Value = -3432.5
# Row 122
$ws.Range("H122").Value = 3546.4849
$ws.Range("I122").Value = 3314.125
$ws.Range("J122").Value = 4166.1113
$ws.Range("K122").Value = 9942.375
$ws.Range("L122").Value = 12498.3339
$ws.Range("M122").Value = -7492.375
$ws.Range("N122").Value = -17398.3339
# Row 136
$ws.Range("H136").Value = 7765.7
$ws.Range("I136").Value = 8690.462
$ws.Range("K136").Value = 26071.386
$ws.Range("M136").Value = -23521.386
